$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report gained a second out-of-stock line item ("LACRITEARS EYE DROPS 15 ML").
# That pushes the existing totals row and the footer/timestamp row down by one row.
# Insert a fresh row at 8 so the old row 8 (totals) becomes row 9 and the old row 9
# (footer) becomes row 10 - mirrors what the generating app does when it re-renders
# the sheet with one more row in the item table.
$ws.Rows("8:8").Insert()

# Seed the new row 8 from row 7's formatting/merges (same visual row "template"),
# then overwrite with the new item's data.
$ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))
$ws.Rows("8:8").RowHeight = 24.75

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "LACRITEARS EYE DROPS 15 ML"
$ws.Range("H8").Value = "1:0"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "49.00"
$ws.Range("P8").Value = "49.0000"
$ws.Range("Q8").Value = "1:0"

# Q8 should carry the same style as Q7 (readingOrder variant of style 8); the
# template copy above leaves it on the generic style, so restore it explicitly.
$ws.Range("Q8").Style = $ws.Range("Q7").Style

# Row 9 (the shifted-down totals row) now sums both item rows' prices and grows
# a touch taller to match the regenerated layout.
$ws.Range("P9").Value = 75.73
$ws.Rows("9:9").RowHeight = 25.5
